$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.050.44"
$ws.Range("E2").Value = "  -4.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.816.32"
$ws.Range("E3").Value = "  -4.53%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.33"
$ws.Range("E5").Value = "  -1.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.61"
$ws.Range("E6").Value = "  +2.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.659"
$ws.Range("E7").Value = "  -3.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.735"
$ws.Range("E9").Value = "  -2.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.37"
$ws.Range("E11").Value = "  -3.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000315"
$ws.Range("E12").Value = "  -1.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.14"
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.415.96"
$ws.Range("E14").Value = "  -4.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.830.88"
$ws.Range("E15").Value = "  -4.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.53"
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.68"
$ws.Range("E17").Value = "  -3.51%  "
$ws.Range("E18").Value = "  -5.99%  "
$ws.Range("E19").Value = "  -2.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.814.80"
$ws.Range("E20").Value = "  -4.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "431.90"
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.68"
$ws.Range("E22").Value = "  -2.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "92.66"
$ws.Range("E23").Value = "  -3.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.24"
$ws.Range("E24").Value = "  -5.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.72"
$ws.Range("E25").Value = "  -3.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.37"
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.95"
$ws.Range("E27").Value = "  -9.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.95"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.35"
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.63"
$ws.Range("E30").Value = "  -5.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.12"
$ws.Range("E31").Value = "  +3.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.33"
$ws.Range("E32").Value = "  -3.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "47.36"
$ws.Range("E33").Value = "  -2.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.124"
$ws.Range("E34").Value = "  -5.40%  "
$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0980"
$ws.Range("E35").Value = "  +6.45%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "68.12"
$ws.Range("E36").Value = "  -3.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "635.75"
$ws.Range("E37").Value = "  -5.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.427"
$ws.Range("E38").Value = "  -2.41%  "
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.144"
$ws.Range("E40").Value = "  -1.38%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.28"
$ws.Range("E42").Value = "  +24.53%  "
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.91"
$ws.Range("E43").Value = "  +11.80%  "
$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.21"
$ws.Range("E44").Value = "  -3.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0464"
$ws.Range("E45").Value = "  -5.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.77"
$ws.Range("E46").Value = "  -8.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.142"
$ws.Range("E47").Value = "  -5.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.825.86"
$ws.Range("E48").Value = "  -1.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.76"
$ws.Range("E49").Value = "  -17.42%  "
$ws.Range("E50").Value = "  -5.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000270"
$ws.Range("E51").Value = "  -1.34%  "
